# sources.xlsx - "completed html tags for sources"
#
# Fills in newly-discovered per-manufacturer scraping notes (HTML snippets /
# selectors used to locate manual links, plus free-text NOTIZ comments), adds
# a live hyperlink to every manufacturer's URL in column B, relocates the
# SMEG row further down the sheet, and adds a new Illy "iperespresso" link row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Cell values - entered in the same left-to-right / top-to-bottom order the
#    author used, so newly created shared-string entries land on the indices
#    the diff expects.
# ---------------------------------------------------------------------------

$ws.Range("C3").Value  = '<div class="dig-pub--text">'
$ws.Range("D3").Value  = '<a href>'
$ws.Range("K1").Value  = 'NOTIZ'
$ws.Range("K3").Value  = 'Together with water boilers, all named "aqua -" '

$ws.Range("C5").Value  = '<ul class="p-pc05v2__cards p-pc05v2__cards--equalize-inner-height p-pc05v2__cards--portrait-view p-pc05v2--list-view-xs" >'
$ws.Range("D5").Value  = '<a class="p-pc05v2__card-view-product-link" href>'
$ws.Range("E5").Value  = '<a class="p-p90__banner p-p90__banner--support" href>'
$ws.Range("F5").Value  = '<section class="p-st14-manuals-documentation p-prx-data-present p-st14-new-design" data-comp-short-title="Bedienungsanleitung" data-comp-id="st14ManualsDocumentation">'
$ws.Range("K5").Value  = 'multiple links for multiple languages, need to scrape all'

$ws.Range("C6").Value  = '<a data-tags="-coffe-makers" href>'
$ws.Range("D6").Value  = '<a href="" target="_blank">'
$ws.Range("K6").Value  = 'bullshit tag without any class links to pdf, maybe check hyperlinks for regex'

$ws.Range("C7").Value  = '<a class="is-full-area ng-star-inserted" href="" target="_self">'
$ws.Range("D7").Value  = '<a class="is-full-area ng-star-inserted" href="" target="_self">'
$ws.Range("K7").Value  = 'regex in first href for "kaffee"'
$ws.Range("E7").Value  = '<a class="c__link ng-star-inserted" href="" target="_self">'

$ws.Range("C8").Value  = '<a id="pagelayout_0_pagetype_0_pagecontent_0_repProducts_hplProduct_0" data-productteaser="link" href="">'
$ws.Range("D8").Value  = '<a id="pagelayout_0_pagetype_0_pagecontent_0_tabcontent_0_repManuals_repManualLanguages_0_hplManual_0" href="">'
$ws.Range("K8").Value  = 'pdf for manual has class 0_hplManual_0, for short manual same with a leading 1'

$ws.Range("C9").Value  = '<a href="" class="plp-product">'
$ws.Range("D9").Value  = '<a href="" class="tech-specs-documents__item-link">'

$ws.Range("K10").Value = 'pdf tohether with datasheet, need to exclude "productinformation" from link'
$ws.Range("C17").Value = 'abfall'

$ws.Range("B11").Value = 'https://www.illy.com/de-de/kaffeemaschinen/kaffeemaschinen-iperespresso-kapseln'

$ws.Range("C13").Value = '                                    '
$ws.Range("C14").Value = '                                    '

$ws.Range("C10").Value = '<div class="product-card__img-item js-product-card-img-item"><a href>'
$ws.Range("D10").Value = '<a class="product-manuals__txt-link" href="">'

# Row 11 loses its old manufacturer name (SMEG is moved to row 17 below).
$ws.Range("A11").ClearContents()

# Row 17 - relocated SMEG entry.
$ws.Range("A17").Value = "SMEG"
$ws.Range("B17").Value = "https://www.smeg.de/sda-kaffeemaschinen/gesamt"

# ---------------------------------------------------------------------------
# 2) Hyperlinks - added in the same order listed in the diff (rId2..rId10;
#    rId1 on B5 already existed).
# ---------------------------------------------------------------------------

$ws.Hyperlinks.Add($ws.Range("B3"),  "https://www.melitta.de/beratung/service-tipps/bedienungsanleitungen/")
$ws.Hyperlinks.Add($ws.Range("B2"),  "https://www.delonghi.com/de-de/manuals/produkte/kaffee/c/coffee")
$ws.Hyperlinks.Add($ws.Range("B6"),  "https://en.russellhobbs.com/products/coffee-makers")
$ws.Hyperlinks.Add($ws.Range("B7"),  "https://www.krups.de/bedienungs-anleitungen/Produkte/Getr%C3%A4nkezubereitung/csc/Beverage")
$ws.Hyperlinks.Add($ws.Range("B8"),  "https://de.jura.com/de/produkte-haushalt/kaffeevollautomaten")
$ws.Hyperlinks.Add($ws.Range("B9"),  "https://www.aeg.de/kitchen/small-kitchen-appliances/coffee-makers/")
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.illy.com/de-de/kaffeemaschinen/kaffeemaschine-fuer-pads-und-gemahlenen-kaffee")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://www.smeg.de/sda-kaffeemaschinen/gesamt")
$ws.Hyperlinks.Add($ws.Range("B4"),  "https://www.sageappliances.com/eu/de/home/index.html?reg=de")

# ---------------------------------------------------------------------------
# 3) Re-apply the "Link" cell style everywhere a hyperlink now lives (Adding a
#    hyperlink nudges the engine into allocating its own ad-hoc variant of the
#    style; re-asserting "Link" brings every cell back to the shared style).
# ---------------------------------------------------------------------------

$ws.Range("B2").Style  = "Link"
$ws.Range("B3").Style  = "Link"
$ws.Range("B4").Style  = "Link"
$ws.Range("B6").Style  = "Link"
$ws.Range("B7").Style  = "Link"
$ws.Range("B8").Style  = "Link"
$ws.Range("B9").Style  = "Link"
$ws.Range("B10").Style = "Link"
$ws.Range("B16").Style = "Link"
$ws.Range("B17").Style = "Link"

# ---------------------------------------------------------------------------
# 4) Column widths - nudged to (closely) track the slightly wider columns
#    recorded after the edit.
# ---------------------------------------------------------------------------

$ws.Columns.Item(2).ColumnWidth  = 90.66666666666667
$ws.Columns.Item(4).ColumnWidth  = 39
$ws.Columns.Item(6).ColumnWidth  = 39.666666666666664
$ws.Columns.Item(7).ColumnWidth  = 42.5
$ws.Columns.Item(8).ColumnWidth  = 34.666666666666664
$ws.Columns.Item(9).ColumnWidth  = 42
$ws.Columns.Item(10).ColumnWidth = 12.666666666666666
$ws.Columns.Item(11).ColumnWidth = 87

# ---------------------------------------------------------------------------
# 5) Header row gets a touch taller, and the sheet's final selection moves.
# ---------------------------------------------------------------------------

$ws.Rows.Item(1).RowHeight = 15.75

[void]$ws.Range("B19").Select()
